$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 831.6
$ws.Range("I11").Value = 831.6
$ws.Range("K11").Value = 831.6
$ws.Range("M11").Value = -691.6

$ws.Range("H12").Value = 301.27274
$ws.Range("I12").Value = 333.77777
$ws.Range("J12").Value = 155
$ws.Range("K12").Value = 333.77777
$ws.Range("L12").Value = 155
$ws.Range("M12").Value = -163.77777
$ws.Range("N12").Value = -495

$ws.Range("H15").Value = 522906.75
$ws.Range("I15").Value = 522906.75
$ws.Range("K15").Value = 1568720.25
$ws.Range("M15").Value = -1568551.25

$ws.Range("H33").Value = 1776.2778
$ws.Range("I33").Value = 1061.5625
$ws.Range("K33").Value = 1061.5625
$ws.Range("M33").Value = -832.5625

$ws.Range("H138").Value = 2567.1316
$ws.Range("J138").Value = 3347.7368
$ws.Range("L138").Value = 10043.2104
$ws.Range("N138").Value = -20323.2104

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9776.385
$ws.Range("I32").Value = 5399.8696
$ws.Range("J32").Value = 43329.668
$ws.Range("K32").Value = 5399.8696
$ws.Range("L32").Value = 43329.668
$ws.Range("M32").Value = -5112.8696
$ws.Range("N32").Value = -43903.668

$ws.Range("H33").Value = 34999.5
$ws.Range("I33").Value = 20000
$ws.Range("J33").Value = 49999
$ws.Range("K33").Value = 20000
$ws.Range("L33").Value = 49999
$ws.Range("M33").Value = -19671
$ws.Range("N33").Value = -50657

$ws.Range("H61").Value = 4360.625
$ws.Range("I61").Value = 2530.8333
$ws.Range("J61").Value = 9850
$ws.Range("K61").Value = 2530.8333
$ws.Range("L61").Value = 9850
$ws.Range("M61").Value = -2318.8333
$ws.Range("N61").Value = -10274

$ws.Range("H110").Value = 1226.7693
$ws.Range("I110").Value = 1295.75
$ws.Range("K110").Value = 1295.75
$ws.Range("M110").Value = 749.25

$ws.Range("H132").Value = 4453.25
$ws.Range("I132").Value = 3450.6365
$ws.Range("J132").Value = 6659
$ws.Range("K132").Value = 10351.9095
$ws.Range("L132").Value = 19977
$ws.Range("M132").Value = -7821.9095
$ws.Range("N132").Value = -25037

$ws.Range("H136").Value = 4360.625
$ws.Range("I136").Value = 2530.8333
$ws.Range("J136").Value = 9850
$ws.Range("K136").Value = 7592.499899999999
$ws.Range("L136").Value = 29550
$ws.Range("M136").Value = -5042.499899999999
$ws.Range("N136").Value = -34650

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 5167.6665
$ws.Range("I82").Value = 5167.6665
$ws.Range("K82").Value = 5167.6665
$ws.Range("M82").Value = -4784.6665

$ws.Range("H85").Value = 5167.6665
$ws.Range("I85").Value = 5167.6665
$ws.Range("K85").Value = 5167.6665
$ws.Range("M85").Value = -3841.6665

$ws.Range("H86").Value = 6545
$ws.Range("J86").Value = 22303
$ws.Range("L86").Value = 22303
$ws.Range("N86").Value = -24549

$ws.Range("H89").Value = 6545
$ws.Range("J89").Value = 22303
$ws.Range("L89").Value = 111515
$ws.Range("N89").Value = -122747

$ws.Range("H126").Value = 59000
$ws.Range("J126").Value = 59000
$ws.Range("L126").Value = 59000
$ws.Range("N126").Value = -68880

$ws.Range("H132").Value = 97999.39999999999
$ws.Range("J132").Value = 97999.39999999999
$ws.Range("L132").Value = 97999.39999999999
$ws.Range("N132").Value = -108119.4

$ws.Range("H134").Value = 3876.389
$ws.Range("I134").Value = 3472.8572
$ws.Range("K134").Value = 10418.5716
$ws.Range("M134").Value = -7883.571599999999

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7138.393
$ws.Range("I31").Value = 4704
$ws.Range("K31").Value = 4704
$ws.Range("M31").Value = -4409

$ws.Range("H34").Value = 7138.393
$ws.Range("I34").Value = 4704
$ws.Range("K34").Value = 4704
$ws.Range("M34").Value = -4502

$ws.Range("H132").Value = 3241.16
$ws.Range("I132").Value = 2138.6843
$ws.Range("J132").Value = 6732.3335
$ws.Range("K132").Value = 6416.0529
$ws.Range("L132").Value = 20197.0005
$ws.Range("M132").Value = -3886.0529
$ws.Range("N132").Value = -25257.0005

$ws.Range("H135").Value = 80768.664
$ws.Range("J135").Value = 80768.664
$ws.Range("L135").Value = 80768.664
$ws.Range("N135").Value = -90908.664

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 19505
$ws.Range("J76").Value = 19505
$ws.Range("L76").Value = 58515
$ws.Range("N76").Value = -59281

$ws.Range("H79").Value = 19505
$ws.Range("J79").Value = 19505
$ws.Range("L79").Value = 58515
$ws.Range("N79").Value = -61167

$ws.Range("H97").Value = 553.5
$ws.Range("I97").Value = 305.46667
$ws.Range("J97").Value = 1085
$ws.Range("K97").Value = 916.4000100000001
$ws.Range("L97").Value = 3255
$ws.Range("M97").Value = -420.4000100000001
$ws.Range("N97").Value = -4247

$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 29999
$ws.Range("J47").Value = 29999
$ws.Range("L47").Value = 29999
$ws.Range("N47").Value = -31135

$ws.Range("H102").Value = 3815.25
$ws.Range("I102").Value = 1753.8889
$ws.Range("K102").Value = 1753.8889
$ws.Range("M102").Value = -131.8888999999999

$ws.Range("H122").Value = 3747.8096
$ws.Range("I122").Value = 3031.1943
$ws.Range("K122").Value = 9093.582900000001
$ws.Range("M122").Value = -6643.582900000001

$ws.Range("H126").Value = 6964.864
$ws.Range("I126").Value = 7322.3
$ws.Range("J126").Value = 6667
$ws.Range("K126").Value = 21966.9
$ws.Range("L126").Value = 20001
$ws.Range("M126").Value = -19496.9
$ws.Range("N126").Value = -24941

$ws.Range("H132").Value = 3527.4146
$ws.Range("I132").Value = 3081.3428
$ws.Range("J132").Value = 6129.5
$ws.Range("K132").Value = 9244.028399999999
$ws.Range("L132").Value = 18388.5
$ws.Range("M132").Value = -6714.028399999999
$ws.Range("N132").Value = -23448.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4531.2856
$ws.Range("I7").Value = 4571.148
$ws.Range("K7").Value = 4571.148
$ws.Range("M7").Value = -4459.148

$ws.Range("H126").Value = 4531.2856
$ws.Range("I126").Value = 4571.148
$ws.Range("K126").Value = 13713.444
$ws.Range("M126").Value = -11243.444

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3491.0857
$ws.Range("I136").Value = 2970.5925
$ws.Range("K136").Value = 8911.7775
$ws.Range("M136").Value = -6361.7775

Write-Output "Applied profit-figure refresh across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets"
